# Updates the cryptos list table (Price / Volume(1h) columns, and the
# NEARProtocol / EnergySwap row swap) to match the new scrape snapshot.
#
# Note: several "Price" values look like plain numbers (e.g. "1.001"),
# but in the source data they are text strings (some even have
# significant trailing zeros, e.g. "1.000"). Assigning such a string
# directly to Range.Value lets Excel auto-convert it to a float and
# lose the formatting/precision, so for those cells we briefly force
# a Text number format, assign the literal string, then restore the
# default "Normal" style so the cell's style index is unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 44 & 45 swap places (NEARProtocol <-> EnergySwap) with new values ---
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.05%  "

$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.198"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.73%  "

# --- Price (D) / Volume(1h) (E) updates for the remaining rows ---
$ws.Range("D2").Value = "30.046.48"
$ws.Range("E2").Value = "  -0.48%  "

$ws.Range("D3").Value = "1.914.94"
$ws.Range("E3").Value = "  +0.51%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5043"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4018"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08248"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.58%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.107"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.05"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.03"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.41%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.426"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.23%  "

$ws.Range("D14").Value = "1.914.50"
$ws.Range("E14").Value = "  +0.28%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.286"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.94%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.83%  "

$ws.Range("E18").Value = "  -1.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06502"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.67%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9997"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.947"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.88%  "

$ws.Range("D23").Value = "30.081.98"
$ws.Range("E23").Value = "  -0.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.199"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "22.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.72%  "

$ws.Range("D27").Value = "2.135.48"
$ws.Range("E27").Value = "  +0.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.283"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.126"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.52%  "

$ws.Range("E32").Value = "  -2.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.995"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.775"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02443"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.347"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.92%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06416"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2162"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.67%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6548"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.59%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.751"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.198"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.38%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.16%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.221"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6033"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.639"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.214"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "78.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.73%  "

$ws.Range("E51").Value = "  -2.97%  "
